$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first match row (old row 2: Velez Sarsfield - Tigre); remaining rows shift up
$ws.Rows.Item(2).Delete()

# Row 2 now holds the former row 3 data (Pereira - Atl. Nacional); refresh a few odds
$ws.Range("G2").Value = 2.57
$ws.Range("I2").Value = 2.85
$ws.Range("M2").Value = 1.11
$ws.Range("N2").Value = 6.5

# Row 3 now holds the former row 4 data (Tepatitlan de Morelos - Tapatio); refresh odds
$ws.Range("H3").Value = 2.85
$ws.Range("I3").Value = 2.55
$ws.Range("M3").Value = 1.08
$ws.Range("N3").Value = 7.5
$ws.Range("O3").Value = 1.37
$ws.Range("P3").Value = 2.62
$ws.Range("Q3").Value = 2.07
$ws.Range("R3").Value = 1.6
$ws.Range("V3").Value = 1.85
$ws.Range("W3").Value = 8.25
$ws.Range("X3").Value = 14.5
$ws.Range("Y3").Value = 10.25
$ws.Range("AA3").Value = 26
$ws.Range("AB3").Value = 35
$ws.Range("AC3").Value = 7.5
$ws.Range("AF3").Value = 70
$ws.Range("AG3").Value = 600
$ws.Range("AH3").Value = 7.3
$ws.Range("AI3").Value = 12.5
$ws.Range("AK3").Value = 30
